# The underlying source data for several field-observation records was
# re-synced, which changed the row order in which the records appear
# (the record content moved between row positions while the sheet's
# layout/header stayed the same). This script reproduces that reshuffle
# by moving the full row contents (columns A:AY) between the affected
# row numbers:
#   row 8  <- old row 9,   row 9  <- old row 8             (swap)
#   row 15 <- old row 16,  row 16 <- old row 15             (swap)
#   row 19 <- old row 21,  row 20 <- old row 19,  row 21 <- old row 20  (rotate)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteValues = -4163
$lastCol = "AY"
$scratchBase = 1000

function Copy-RowValues($srcRow, $dstRow) {
    $srcAddr = "A" + $srcRow + ":" + $lastCol + $srcRow
    $dstAddr = "A" + $dstRow + ":" + $lastCol + $dstRow
    $src = $ws.Range($srcAddr)
    $dst = $ws.Range($dstAddr)
    $src.Copy()
    $dst.PasteSpecial($xlPasteValues)
}

function Clear-RowContents($row) {
    $addr = "A" + $row + ":" + $lastCol + $row
    $ws.Range($addr).ClearContents()
}

function Invoke-RowCycle([int[]]$rows) {
    # new content of $rows[i] <- old content of $rows[i-1] (wrapping around)
    $n = $rows.Length

    # Stash every row's original content into scratch rows first, since the
    # destinations overlap with the sources.
    for ($i = 0; $i -lt $n; $i++) {
        $srcRow = $rows[$i]
        $scratchRow = $scratchBase + $i
        Copy-RowValues $srcRow $scratchRow
    }

    # Clear the originals so stale cells that have no counterpart in the
    # incoming row don't linger behind.
    for ($i = 0; $i -lt $n; $i++) {
        $row = $rows[$i]
        Clear-RowContents $row
    }

    # Write each row's new content: row[i] <- old row[i-1] (wrapping)
    for ($i = 0; $i -lt $n; $i++) {
        $prev = $i - 1
        if ($prev -lt 0) { $prev = $n - 1 }
        $scratchRow = $scratchBase + $prev
        $dstRow = $rows[$i]
        Copy-RowValues $scratchRow $dstRow
    }

    # Clean up scratch rows.
    for ($i = 0; $i -lt $n; $i++) {
        $scratchRow = $scratchBase + $i
        Clear-RowContents $scratchRow
    }
}

# row 8 <- old 9, row 9 <- old 8
$group1 = @(8, 9)
Invoke-RowCycle $group1

# row 15 <- old 16, row 16 <- old 15
$group2 = @(15, 16)
Invoke-RowCycle $group2

# row 20 <- old 19, row 21 <- old 20, row 19 <- old 21
$group3 = @(20, 21, 19)
Invoke-RowCycle $group3

Write-Host "Row reshuffle applied."
